$d = $word.ActiveDocument

# The "_GoBack" bookmark sits between two runs in the Game Controls
# paragraph: "...SPACEBAR jumps. C" [bookmark] " crouches. SHIFT walks. "
# We need the run before the bookmark to absorb the " crouches. SHIFT
# walks. " text, and the (now stale) run after the bookmark to instead
# read the new "ESC to return to menu from in game." sentence.
$bm = $word.ActiveDocument.Bookmarks("_GoBack")

# Range that currently holds " crouches. SHIFT walks. " (right after the
# bookmark). Replace it first since it is further along in the document.
$after = $d.Range($bm.End, $bm.End + 24)
$after.Text = "ESC to return to menu from in game."

# Range that currently holds "SPACEBAR jumps. C" (right before the
# bookmark). Append the crouch/walk text that used to live after the
# bookmark.
$before = $d.Range($bm.Start - 17, $bm.Start)
$before.Text = "SPACEBAR jumps. C crouches. SHIFT walks. "
